$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testCitizen")
$ws.Activate()

# Row 1, column B: "ubs111" -> "ulis111" (duplicate of A1's value)
$ws.Range("B1").Value = "ulis111"

# Update the selected/active cell on the sheet
$ws.Range("N10").Select()
